# Generate Report for Handback
# - Update the "Ready for handoff" status to "Handback transform failed" for the
#   2d619095-8c47-40a5-ba95-de4db1e38de8 file (Overview sheet + per-language sheets).
# - Record an Error Detail message on the zh-cn and de-de sheets explaining the
#   handback filename mismatch.
# - Widen the "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 2d619095-... entry, columns E (zh-cn) and F (de-de)
$wsOverview.Cells.Item(3, 5).Value = $newStatus
$wsOverview.Cells.Item(3, 6).Value = $newStatus

# zh-cn / de-de sheets: row 3 is the 2d619095-... entry, column C is "Status"
$wsZhCn.Cells.Item(3, 3).Value = $newStatus
$wsDeDe.Cells.Item(3, 3).Value = $newStatus

# Column P ("Error Detail") gets the failure detail message, and is widened.
$wsZhCn.Cells.Item(3, 16).Value = "Handback file name: bsbmzipl.ynp is different with handoff file name: 2d619095-8c47-40a5-ba95-de4db1e38de8.3c175f97ed7582c4fc9b5119441d66466189bdbf.zh-cn."
$wsDeDe.Cells.Item(3, 16).Value = "Handback file name: bsbmzipl.ynp is different with handoff file name: 2d619095-8c47-40a5-ba95-de4db1e38de8.3c175f97ed7582c4fc9b5119441d66466189bdbf.de-de."

$wsZhCn.Columns.Item(16).ColumnWidth = 39.2
$wsDeDe.Columns.Item(16).ColumnWidth = 39.2
